# edit.ps1 -- applies the LOM3226.docx content update described by the diff.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Text blocks (single-quoted here-strings -> no variable interpolation)
# ---------------------------------------------------------------------
$objPtAnchor = @'
Apresentar o formalismo para descrição de sistemas quânticos. Estudar diversas aplicações da equação de Schroedinger independente do tempo. Descrever a estrutura eletrônica de átomos e moléculas.
'@
$objEn = @'
To present the formalism for the description of quantum systems. Study several applications of the time-independent Schroedinger equation. Describe the electronic structure of atoms and molecules.
'@
$summaryPtAnchor = @'
Introdução aos conceitos da Mecânica Quântica. • Ferramentas matemáticas da Mecânica Quântica. A equação de Schroedinger e aplicações unidimensionais e tridimensionais. Problemas em coordenadas retangulares. Problemas em coordenadas esféricas. Átomos com um elétron. Teoria geral. • Propriedades gerais do momento angular.
'@
$summaryEn = @'
• Introduction to the concepts of Quantum Mechanics. • Mathematical tools of Quantum Mechanics. • The Schrödinger equation and one- and three-dimensional applications. • Quantum formalism. • Problems in rectangular coordinates and spherical coordinates. • Hydrogen atoms and orbitals. • General properties of angular momentum. • Spin. • Fermions and bosons.
'@
$programaOldPt = @'
Origens das ideias fundamentais da Mecânica Quântica.Dualidade onda partícula. Principio de Heisenberg.• Os postulados da Mecânica Quântica. Ferramentas matemáticas da Mecânica Quântica. O Espaço de Hilbert e a Equação de Onda. Notação de Dirac. Operadores e Bases. Representação matricial. A equação de Schroedinger e aplicações unidimensionais. Barreira de potencial. Poço de potencial. Oscilador harmônico. Problemas tridimensionais. Problemas em coordenadas retangulares. Problemas em coordenadas esféricas. Átomos com um elétron. Teoria geral.• Momento angular. Spin do elétron.• Propriedades gerais do momento angular
'@
$programaNewPt = @'
• Equação de Schrödinger. • Função de onda e interpretação estatística da mecânica quântica. • Valores esperados e operadores. Os operadores posição e momento; operadores energia cinética e potencial; o operador Hamiltoniano. • A equação de Schrödinger independente do tempo. Separação de variáveis e estados estacionários. • Aplicações unidimensionais:  poço quadrado infinito; oscilador harmônico; partícula livre;  transformada de Fourier e sua relação com o princípio da incerteza de Heisenberg; Poços e barreiras de potencial. • Formalismo quântico: opserváveis e operadores hermitianos. Estados determinados, autoestados e autovalores de operadores hermitianos. Base de autoestados; interpretação estatística generalizada: medidas de observáveis e suas probabilidades. Comutadores e operadores que compartilham autoestados ; princípio da incerteza generalizado. • Mecânica Quântica em três dimensões. • Átomo de hidrogênio: modelo de Bohr e o número quântico principal. Solução completa e os demais números quânticos. • Coordenadas esféricas e Momento angular.  • Momento angulas de spin. • Problemas de muitos corpos. • Partículas idênticas: férmions e bósons.
'@
$programaNewEn = @'
• Schrödinger's equation. • Wave function and statistical interpretation of quantum mechanics. • Expected values and operators. The position and moment operators; kinetic and potential energy operators; the Hamiltonian operator. • The time-independent Schrödinger equation. Separation of variables and steady states. • One-dimensional applications: infinite square well; harmonic oscillator; free particle; Fourier transform and its relationship with the Heisenberg uncertainty principle; Potential square wells and barriers. • Quantum formalism: hermitian operators and observables. Determined states, eigenstates and eigenvalues of Hermitian operators. Basis of Eigenstates; generalized statistical interpretation: measures of observables and their probabilities. Comutators and operators that share eigenstates; generalized uncertainty principle. • Quantum Mechanics in three dimensions. • Hydrogen atom: Bohr model and the principal quantum number. Complete solution and the other quantum numbers. • Spherical coordinates and Angular momentum. • Spin angular momentum. • Many-body problems. • Identical particles: fermions and bosons.
'@
$criterioOld = @'
Média aritmética de duas provas sendo a primeira com peso 1 e a segunda com peso 2.
'@
$criterioNew = @'
Média aritmética de três provas: P1 (peso 1), P2 (peso 1) e P3 (peso 2).
'@
$bulletCarlos = @'
6279110 - Carlos Alberto Moreira dos Santos
'@
$bulletDurval = @'
6495737 - Durval Rodrigues Junior
'@

# ---------------------------------------------------------------------
# 1) Activation date: 01/01/2020 -> 01/01/2023
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ativação: 01/01/2023", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert italic English translation paragraph after the PT 'Objetivos' text
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute($objPtAnchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter($objEn)
$rng.Font.Italic = 1

# ---------------------------------------------------------------------
# 3) Remove the Carlos Alberto Moreira dos Santos and Durval Rodrigues
#    Junior lines (with their trailing line breaks) from Docente(s)
# ---------------------------------------------------------------------
$rng = $d.Content
$delText = $bulletCarlos + [char]11 + $bulletDurval + [char]11
$rng.Find.Execute($delText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Delete()

# ---------------------------------------------------------------------
# 4) Insert italic English translation paragraph after the PT 'Programa
#    resumido' text
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute($summaryPtAnchor, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter($summaryEn)
$rng.Font.Italic = 1

# ---------------------------------------------------------------------
# 5) Replace the PT 'Programa' body text, then insert the new italic EN
#    translation paragraph right after it
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute($programaOldPt, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = $programaNewPt
$rng.Collapse(0)
$rng.InsertParagraphAfter() | Out-Null
$rng.Move(1, 1) | Out-Null
$rng.InsertAfter($programaNewEn)
$rng.Font.Italic = 1

# ---------------------------------------------------------------------
# 6) Avaliação criterion: two-exam average -> three-exam average
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute($criterioOld, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Text = $criterioNew
